function HexToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = HexToOleColor("000000")
$cs.Item(2).RGB  = HexToOleColor("FFFFFF")
$cs.Item(3).RGB  = HexToOleColor("44546A")
$cs.Item(4).RGB  = HexToOleColor("E7E6E6")
$cs.Item(5).RGB  = HexToOleColor("5B9BD5")
$cs.Item(6).RGB  = HexToOleColor("ED7D31")
$cs.Item(7).RGB  = HexToOleColor("A5A5A5")
$cs.Item(8).RGB  = HexToOleColor("FFC000")
$cs.Item(9).RGB  = HexToOleColor("4472C4")
$cs.Item(10).RGB = HexToOleColor("70AD47")
$cs.Item(11).RGB = HexToOleColor("0563C1")
$cs.Item(12).RGB = HexToOleColor("954F72")
